# Automatic update of files.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C from 45593 to 45594
# for every data row (2 through 31).
$ws.Range("C2:C31").Value2 = 45594

# Rows 30 and 31 had their beteckning (A) and area (G) values swapped.
$ws.Range("A30").Value2 = "A 46085-2024"
$ws.Range("G30").Value2 = 1.9

$ws.Range("A31").Value2 = "A 46082-2024"
$ws.Range("G31").Value2 = 1
